# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect newly scraped totals.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 686
$ws1.Range("F3").Value = 30
$ws1.Range("F4").Value = 231
$ws1.Range("F5").Value = 2025
$ws1.Range("F7").Value = 3359
$ws1.Range("F9").Value = 818

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 686
$ws4.Range("F3").Value = 30
$ws4.Range("F5").Value = 231
$ws4.Range("F6").Value = 2025
$ws4.Range("F8").Value = 3359
$ws4.Range("F10").Value = 818
